$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Date - 22.01.2024"

$text_D4 = @"
DEL.N1LN3L.WD.D.GPD.DMP : last exported on 20.01.24
DEL.N1LN3L.WD.D.GPL.DMP : last exported on 20.01.24
DEL.N1LN3L.WD.D.GTR.DMP : last exported on 21.01.24
DEL.N1LKQQ.WD.D.UPD.DMP : last exported on 20.01.24
DEL.N1LN3L.WD.D.ASPD.DMP : last exported on 20.01.24
DEL.N1LN3L.WD.D.ASPL.DMP : last exported on 20.01.24
DEL.N1LN0R.WD.D.CPD.DMP : last exported on 20.01.24
DEL.N1LN0R.WD.D.CPL.DMP : last exported on 20.01.24
DEL.KT6N5M.SL.D.GGO.ZIP : last exported on 19.01.24
DEL.KT6N5M.SL.D.UGO.ZIP : last exported on 19.01.24
DEL.KT6N5M.SL.D.ASGO.ZIP : last exported on 19.01.24
"@
$ws.Range("D4").Value = $text_D4

$text_D6 = @"
DEL.N1LN3L.DN.F.GPD.DMP : last exported on 22.01.24
DEL.N1LN3L.DN.F.GPL.DMP : last exported on 22.01.24
DEL.N1LN3L.DN.F.GTR.DMP : last exported on 22.01.24
DEL.N1LKQQ.DN.F.UPD.DMP : last exported on 22.01.24
DEL.N1LN3L.DN.F.ASPD.DMP : last exported on 22.01.24
DEL.N1LN3L.DN.F.ASPL.DMP : last exported on 22.01.24
DEL.N1LN0R.DN.F.CPD.DMP : last exported on 22.01.24
DEL.N1LN0R.DN.F.CPL.DMP : last exported on 22.01.24
DEL.KT6E35.SN.F.GGO.ZIP : last exported on 15.01.24
DEL.KT6N5M.SN.F.GGO.ZIP : last exported on 15.01.24
"@
$ws.Range("D6").Value = $text_D6

$text_D7 = @"
DEL.N1LN3L.WI.D.GPD.DMP : last exported on 22.01.24
DEL.N1LN3L.WI.D.GPL.DMP : last exported on 22.01.24
DEL.N1LN3L.WI.D.GTR.DMP : last exported on 22.01.24
DEL.N1LN3L.WI.D.ASPD.DMP : last exported on 22.01.24
DEL.N1LN3L.WI.D.ASPL.DMP : last exported on 22.01.24
DEL.N1LN0R.WI.D.CPD.DMP : last exported on 22.01.24
DEL.N1LN0R.WI.D.CPL.DMP : last exported on 22.01.24
DEL.KT6E35.WI.GGO.ZIP : last exported on 18.01.24
DEL.KT6E35.WI.IMG.GGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.WI.GGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.WI.IMG.GGO.ZIP : last exported on 18.01.24
"@
$ws.Range("D7").Value = $text_D7

$text_D10 = @"
DEL.N1LN3L.TP.D.GPD.DMP : last exported on 21.01.24
DEL.N1LN3L.TP.D.GPL.DMP : last exported on 21.01.24
DEL.N1LN3L.TP.D.GTR.DMP : last exported on 21.01.24
DEL.N1LKQQ.TP.D.UPD.DMP : last exported on 21.01.24
DEL.N1LN3L.TP.D.ASPD.DMP : last exported on 21.01.24
DEL.N1LN3L.TP.D.ASPL.DMP : last exported on 21.01.24
DEL.N1LN0R.TP.D.CPD.DMP : last exported on 21.01.24
DEL.N1LN0R.TP.D.CPL.DMP : last exported on 21.01.24
DEL.KT6E35.TP.D.GGO.ZIP : last exported on 18.01.24
DEL.KT6E35.TP.D.UGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.GGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.UGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.ASGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.UOC.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.ASOC.ZIP : last exported on 18.01.24
DEL.KT6N5M.TP.D.GOC.ZIP : last exported on 18.01.24
DED.KT6I1M.HST2.D.GGO.ZIP : last exported on 18.01.24
"@
$ws.Range("D10").Value = $text_D10

$text_D11 = @"
DEL.N1LN3L.MC.F.GPD.DMP : last exported on 30.11.23
DEL.N1LN3L.MC.F.GPL.DMP : last exported on 30.11.23
DEL.N1LN3L.MC.F.GTR.DMP : last exported on 21.01.24
DEL.N1LKQQ.MC.F.UPD.DMP : last exported on 30.11.23
DEL.N1LN3L.MC.F.ASPD.DMP : last exported on 30.11.23
DEL.N1LN3L.MC.F.ASPL.DMP : last exported on 30.11.23
DGI.KT6R11.MANDANT.TXT : last exported on 30.11.23
"@
$ws.Range("D11").Value = $text_D11

$text_D13 = @"
DEL.N1LN3L.AU.D.GPD.DMP : last exported on 23.12.23
DEL.N1LN3L.AU.D.GPL.DMP : last exported on 23.12.23
DEL.N1LN3L.AU.D.GTR.DMP : last exported on 21.01.24
DEL.N1LN0R.AU.D.CPD.DMP : last exported on 23.12.23
DEL.N1LN0R.AU.D.CPL.DMP : last exported on 23.12.23
DEL.KT6E35.AU.GGO.ZIP : last exported on 23.12.23
DEL.KT6E35.AU.IMG.GGO.ZIP : last exported on 23.12.23
"@
$ws.Range("D13").Value = $text_D13

$text_D15 = @"
DGI.N1LR11.FB.ORGUNIT.CSV : last exported on 21.01.24
DGI.N1LR11.FB.IMPTR.CSV : last exported on 21.01.24
"@
$ws.Range("D15").Value = $text_D15

$text_D16 = @"
DEL.N1LN3L.DP.F.GPD.DMP : last exported on 19.01.24
DEL.N1LN3L.DP.F.GPL.DMP : last exported on 19.01.24
DEL.N1LN3L.DP.F.GTR.DMP : last exported on 21.01.24
DEL.N1LKQQ.DP.F.UPD.DMP : last exported on 19.01.24
DEL.N1LN3L.DP.F.ASPD.DMP : last exported on 19.01.24
DEL.N1LN3L.DP.F.ASPL.DMP : last exported on 19.01.24
DEL.N1LN0R.DP.F.CPD.DMP : last exported on 19.01.24
DEL.N1LN0R.DP.F.CPL.DMP : last exported on 19.01.24
DEL.N1LN3L.DISSORG.ZIP : last exported on 19.01.24
DEL.N1LN3L.DISSAORG.ZIP : last exported on 19.01.24
"@
$ws.Range("D16").Value = $text_D16

$ws.Range("D18").Value = "DEL.N1LE35.TRANS.ZIP : last exported on 21.01.24"

$text_D19 = @"
DEL.N1LN3L.MT.F.GPD.DMP : last exported on 20.01.24
DEL.N1LN3L.MT.F.GPL.DMP : last exported on 20.01.24
DEL.N1LN3L.MT.F.GTR.DMP : last exported on 21.01.24
DEL.N1LKQQ.MT.F.UPD.DMP : last exported on 20.01.24
DEL.N1LN3L.MT.F.ASPD.DMP : last exported on 20.01.24
DEL.N1LN3L.MT.F.ASPL.DMP : last exported on 20.01.24
DEL.N1LN0R.MT.F.CPD.DMP : last exported on 20.01.24
DEL.N1LN0R.MT.F.CPL.DMP : last exported on 20.01.24
DEL.N1LR31.MT.F.GPD.DMP : last exported on 20.01.24
DEL.N1LI4X.MT.F.GPD.ZIP : last exported on 20.01.24
DEL.N1LI4Y.MT.F.GPL.ZIP : last exported on 20.01.24
DEL.N1LGN1.NL.F.GPD.ZIP : last exported on 20.01.24
DEL.N1LKT6.MT.??????.G : last exported on 20.01.24
DEL.N1LKT6.MT.??????.U : last exported on 20.01.24
DEL.N1LKT6.MT.??????.AS : last exported on 20.01.24
DVL.N1LN5X.VLM.DSP : last exported on 20.01.24
DVL.N1LN5X.VLM.WT : last exported on 20.01.24
DDS.N1LR11.DSP : last exported on 20.01.24
DDS.N1LR11.WT : last exported on 20.01.24
DEL.KT6N5M.MT.F.GGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.MT.F.UGO.ZIP : last exported on 18.01.24
DEL.KT6N5M.MT.F.ASGO.ZIP : last exported on 18.01.24
DEL.KT6E35.MT.F.GGO.ZIP : last exported on 18.01.24
DEL.KT6E35.MT.F.UGO.ZIP : last exported on 18.01.24
DEL.KT6N5X.MT.F.GGO.ZIP : last exported on 18.01.24
DEL.KT6N5X.MT.F.UGO.ZIP : last exported on 18.01.24
DEL.KT6N5X.MT.F.ASGO.ZIP : last exported on 18.01.24
"@
$ws.Range("D19").Value = $text_D19

$text_D20 = @"
DEL.N1LN3L.EL.F.GPD.DMP : last exported on 21.01.24
DEL.N1LN3L.EL.F.GPL.DMP : last exported on 21.01.24
DEL.N1LE35.ELFI.G.DATA : last exported on 21.01.24
DEL.N1LN5M.EL.F.GGO.ZIP : last exported on 21.01.24
DEL.N1LN5X.EL.F.GGO.ZIP : last exported on 21.01.24
DEL.N1LE35.EL.F.GGO.ZIP : last exported on 21.01.24
DEL.N1LKT6.EC.??????.G : last exported on 21.01.24
DWT.N1LN3L.WT.??????.ZIP : last exported on 21.01.24
DEL.N1LKQQ.EL.F.UPD.DMP : last exported on 21.01.24
DEL.N1LE35.ELFI.U.DATA : last exported on 21.01.24
DEL.N1LN5M.EL.F.UGO.ZIP : last exported on 21.01.24
DEL.N1LN5X.EL.F.UGO.ZIP : last exported on 21.01.24
DEL.N1LE35.EL.F.UGO.ZIP : last exported on 21.01.24
DEL.N1LKT6.EC.??????.U : last exported on 21.01.24
DEL.N1LN3L.EL.F.ASPD.DMP : last exported on 21.01.24
DEL.N1LN3L.EL.F.ASPL.DMP : last exported on 21.01.24
DEL.N1LN0R.EL.F.CPD.DMP : last exported on 21.01.24
DEL.N1LN0R.EL.F.CPL.DMP : last exported on 21.01.24
DEL.N1LN5M.EL.F.ASGO.ZIP : last exported on 21.01.24
DEL.N1LN5X.EL.F.ASGO.ZIP : last exported on 21.01.24
DEL.N1LE35.EL.F.ASGO.ZIP : last exported on 21.01.24
DEL.N1LKT6.EC.??????.AS : last exported on 21.01.24
DEL.N1LN3L.EL.F.GTR.DMP : last exported on 21.01.24
DGV.N1LR31.EL.F.HST.ZIP : last exported on 21.01.24
DGV.N1LN3L.EL.F.HST.ZIP : last exported on 21.01.24
DGV.N1LCEN.EL.F.ZIP : last exported on 21.01.24
DGV.N1LN7K.EL.F.ZIP : last exported on 21.01.24
DGV.N1LN5X.EL.F.ZIP : last exported on 21.01.24
DEL.KT6KT6.EL.F.DRE.CSV : last exported on 21.01.24
DEL.KT6N5X.EL.F.AU.ZIP : last exported on 18.01.24
"@
$ws.Range("D20").Value = $text_D20
